$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.440.04'
$ws.Range('E2').Value = '  +4.23%  '
$ws.Range('D3').Value = '1.803.79'
$ws.Range('E3').Value = '  +1.63%  '
$ws.Range('D4').Value = "'0.9996"
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Value = "'315.12"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.50%  '
$ws.Range('D6').Value = "'0.9994"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.13%  '
$ws.Range('D7').Value = "'0.5502"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.90%  '
$ws.Range('D8').Value = "'0.3863"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +6.04%  '
$ws.Range('D9').Value = "'0.07606"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.27%  '
$ws.Range('D10').Value = "'42.57"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.27%  '
$ws.Range('E11').Value = '  +3.38%  '
$ws.Range('D12').Value = "'0.9997"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.13%  '
$ws.Range('D13').Value = "'21.23"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.64%  '
$ws.Range('D14').Value = "'6.194"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.20%  '
$ws.Range('D15').Value = "'7.461"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +7.31%  '
$ws.Range('D16').Value = '1.806.48'
$ws.Range('E16').Value = '  +2.19%  '
$ws.Range('D17').Value = "'92.04"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.77%  '
$ws.Range('E18').Value = '  +2.69%  '
$ws.Range('D19').Value = "'0.06442"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.28%  '
$ws.Range('D20').Value = "'0.9993"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.12%  '
$ws.Range('E21').Value = '  +3.66%  '
$ws.Range('E22').Value = '  +2.62%  '
$ws.Range('D23').Value = '28.442.64'
$ws.Range('E23').Value = '  +3.94%  '
$ws.Range('D24').Value = "'11.44"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.50%  '
$ws.Range('D25').Value = "'2.132"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.76%  '
$ws.Range('D26').Value = "'159.08"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.07%  '
$ws.Range('D27').Value = "'20.71"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.05%  '
$ws.Range('D28').Value = "'2.410"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.22%  '
$ws.Range('D29').Value = '2.011.44'
$ws.Range('E29').Value = '  +2.07%  '
$ws.Range('E30').Value = '  +2.49%  '
$ws.Range('E31').Value = '  +6.15%  '
$ws.Range('D32').Value = "'0.1023"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.69%  '
$ws.Range('D33').Value = "'5.772"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.90%  '
$ws.Range('D34').Value = "'3.686"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.88%  '
$ws.Range('D35').Value = "'0.2310"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +14.31%  '
$ws.Range('D36').Value = "'0.06454"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +8.31%  '
$ws.Range('E37').Value = '  +4.48%  '
$ws.Range('D38').Value = "'5.180"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.10%  '
$ws.Range('D39').Value = "'8.824"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +9.46%  '
$ws.Range('D40').Value = "'11.66"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.96%  '
$ws.Range('D41').Value = "'0.6415"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.71%  '
$ws.Range('E42').Value = '  +2.11%  '
$ws.Range('D43').Value = "'0.9991"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.11%  '
$ws.Range('D44').Value = "'1.384"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.41%  '
$ws.Range('D45').Value = "'13.53"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.05%  '
$ws.Range('D46').Value = "'0.5986"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.17%  '
$ws.Range('E47').Value = '  +1.52%  '
$ws.Range('D48').Value = "'127.27"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.33%  '
$ws.Range('D49').Value = "'1.984"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.21%  '
$ws.Range('E50').Value = '  +3.75%  '
$ws.Range('D51').Value = "'0.06894"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.71%  '
